# This script re-applies the "before" values of columns A,B,D,E,F,G,H,Q,R
# (rows 2-14) according to a permutation of rows, producing the "after"
# state described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for the columns that get shuffled across rows 2-14
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$original = @{}
for ($r = 2; $r -le 14; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Mapping: target row -> source row (data pulled from the original source row)
$mapping = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 8
    6  = 2
    7  = 3
    8  = 9
    9  = 10
    10 = 11
    11 = 12
    12 = 13
    13 = 14
    14 = 4
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $srcVals = $original[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $srcVals[$c]
    }
}
